$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 8
$ws.Range("F3").Value = -1
$ws.Range("F4").Value = 1
$ws.Range("F6").Value = -2
$ws.Range("F7").Value = -2
$ws.Range("F10").Value = -2
$ws.Range("F11").Value = 4
$ws.Range("F12").Value = 1
$ws.Range("F13").Value = -2
$ws.Range("F14").Value = 4
$ws.Range("F15").Value = -4
